# Update dashboard, Dispatch, Inventory, History
#
# 1. Sheet2!A1 holds the long Thai "dashboard spec" note. The first bullet
#    ("chong thi 1") needs an extra clause inserted right after
#    "...dtong song" and before " chong thi 2": "thang-mod 500 dto wan".
# 2. The last user selection on Sheet2 moves from F12 to Q6.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")
$cell = $ws2.Range("A1")

$oldText = $cell.Value2
$oldFragment = "ที่ต้องส่ง ช่องที่2"
$newFragment = "ที่ต้องส่งทั้งหมด500ต่อวัน ช่องที่2"

if ($oldText.Contains($oldFragment)) {
    $cell.Value2 = $oldText.Replace($oldFragment, $newFragment)
}

$ws2.Activate() | Out-Null
$ws2.Range("Q6").Select() | Out-Null
